# Regenerate save_data: column G ("K") values are recalculated (std/mean based
# "s_vals" calc), replacing the previous "Strike#" derived values.
# Apply the new computed K values for rows 2-21 in the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 0
    4  = 4
    5  = 1
    6  = 2
    7  = 1
    8  = 5
    9  = 3
    10 = 2
    11 = 4
    12 = 5
    13 = 0
    14 = 3
    15 = 0
    16 = 7
    17 = 5
    18 = 2
    19 = 3
    20 = 1
    21 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
